# Regenerate orders with updated distance/size codes.
# Distances: D51 -> D55, D64 -> D69, D80 -> D86
# Sizes:     S30 -> S31 (S20 / S25 stay the same)
#
# These tokens appear embedded inside many shared strings across several
# columns (Condition, Filename_Left, Filename_Right, Distance, Size), so a
# workbook-wide Find & Replace on the distinct tokens reproduces the rename
# everywhere it occurs (headers like D64/D51/D80 in the Distance column,
# S30 in the Size column, and the composite FaceNN_Dxx_Syy / *_l.png / *_r.png
# strings alike).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters only in that the four tokens are disjoint, so any order is
# safe; replace the distance codes first, then the size code. Results are
# discarded ([void]) so the COM log only reports genuine problems.
[void]$ws.Cells.Replace("D51", "D55")
[void]$ws.Cells.Replace("D64", "D69")
[void]$ws.Cells.Replace("D80", "D86")
[void]$ws.Cells.Replace("S30", "S31")
